$wb = $excel.ActiveWorkbook

# ALC row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 11332.533
$ws.Range("I70").Value = 4995
$ws.Range("J70").Value = 12307.538
$ws.Range("K70").Value = 14985
$ws.Range("L70").Value = 36922.614
$ws.Range("M70").Value = -14715
$ws.Range("N70").Value = -37462.614

# ALC row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 11332.533
$ws.Range("I73").Value = 4995
$ws.Range("J73").Value = 12307.538
$ws.Range("K73").Value = 14985
$ws.Range("L73").Value = 36922.614
$ws.Range("M73").Value = -14049
$ws.Range("N73").Value = -38794.614

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2999
$ws.Range("I113").Value = 2999
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2999
$ws.Range("L113").ClearContents()
$ws.Range("M113").Value = 255
$ws.Range("N113").ClearContents()

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1800.3334
$ws.Range("I137").Value = 1875.5
$ws.Range("J137").Value = 1650
$ws.Range("K137").Value = 5626.5
$ws.Range("L137").Value = 4950
$ws.Range("M137").Value = -3076.5
$ws.Range("N137").Value = -10050

# ARM row 21
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 980
$ws.Range("I21").Value = 666.6667
$ws.Range("J21").Value = 1450
$ws.Range("K21").Value = 666.6667
$ws.Range("L21").Value = 1450
$ws.Range("M21").Value = -292.6667
$ws.Range("N21").Value = -2198

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1103.4783
$ws.Range("I32").Value = 867.7778
$ws.Range("J32").Value = 1952
$ws.Range("K32").Value = 867.7778
$ws.Range("L32").Value = 1952
$ws.Range("M32").Value = -580.7778
$ws.Range("N32").Value = -2526

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2500
$ws.Range("I45").Value = 2500
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2500
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").ClearContents()

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3497.5
$ws.Range("I61").Value = 3497.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3497.5
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 25000
$ws.Range("I74").Value = 25000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 25000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -24126

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 25000
$ws.Range("I77").Value = 25000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 125000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -120632

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3497.5
$ws.Range("I136").Value = 3497.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10492.5
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5260.375
$ws.Range("I105").Value = 5260.375
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 5260.375
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -3513.375

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5699.6665
$ws.Range("I134").Value = 5699.6665
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 17098.9995
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -14563.9995

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1515.25
$ws.Range("I31").Value = 1273.3334
$ws.Range("J31").Value = 1660.4
$ws.Range("K31").Value = 1273.3334
$ws.Range("L31").Value = 1660.4
$ws.Range("M31").Value = -978.3334
$ws.Range("N31").Value = -2250.4

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1515.25
$ws.Range("I34").Value = 1273.3334
$ws.Range("J34").Value = 1660.4
$ws.Range("K34").Value = 1273.3334
$ws.Range("L34").Value = 1660.4
$ws.Range("M34").Value = -1071.3334
$ws.Range("N34").Value = -2064.4

# CRP row 93
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 25000
$ws.Range("I93").Value = 25000
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 25000
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -23128

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 420.1111
$ws.Range("I5").Value = 347.625
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 1042.875
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = -930.875
$ws.Range("N5").Value = -3224

# CUL row 25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()

# CUL row 26
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()

# CUL row 30
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").ClearContents()

# CUL row 94
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 5050
$ws.Range("I94").Value = 5050
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 15150
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -14474

# CUL row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 425.66666
$ws.Range("I121").Value = 463.5
$ws.Range("J121").Value = 350
$ws.Range("K121").Value = 1390.5
$ws.Range("L121").Value = 1050
$ws.Range("M121").Value = -80.5
$ws.Range("N121").Value = -3670

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1411
$ws.Range("I132").Value = 475
$ws.Range("J132").Value = 2035
$ws.Range("K132").Value = 4275
$ws.Range("L132").Value = 18315
$ws.Range("M132").Value = -1745
$ws.Range("N132").Value = -23375

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 420.1111
$ws.Range("I135").Value = 347.625
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 3128.625
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -593.625
$ws.Range("N135").Value = -14070

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 623.3333
$ws.Range("I140").Value = 623.3333
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 1869.9999
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = 3310.0001

# GSM row 53
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 15642.5
$ws.Range("I53").Value = 9285
$ws.Range("J53").Value = 22000
$ws.Range("K53").Value = 9285
$ws.Range("L53").Value = 22000
$ws.Range("M53").Value = -8654
$ws.Range("N53").Value = -23262

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()

# GSM row 82
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()

# GSM row 85
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()

# LTW row 4
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 5180
$ws.Range("I4").Value = 2725
$ws.Range("J4").Value = 15000
$ws.Range("K4").Value = 2725
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = -2612
$ws.Range("N4").Value = -15226

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2647.5789
$ws.Range("I22").Value = 1771.1428
$ws.Range("J22").Value = 3158.8333
$ws.Range("K22").Value = 1771.1428
$ws.Range("L22").Value = 3158.8333
$ws.Range("M22").Value = -1476.1428
$ws.Range("N22").Value = -3748.8333

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2647.5789
$ws.Range("I27").Value = 1771.1428
$ws.Range("J27").Value = 3158.8333
$ws.Range("K27").Value = 1771.1428
$ws.Range("L27").Value = 3158.8333
$ws.Range("M27").Value = -1664.1428
$ws.Range("N27").Value = -3372.8333

# LTW row 28
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H28").Value = 5180
$ws.Range("I28").Value = 2725
$ws.Range("J28").Value = 15000
$ws.Range("K28").Value = 2725
$ws.Range("L28").Value = 15000
$ws.Range("M28").Value = -2493
$ws.Range("N28").Value = -15464

# LTW row 37
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H37").Value = 5180
$ws.Range("I37").Value = 2725
$ws.Range("J37").Value = 15000
$ws.Range("K37").Value = 2725
$ws.Range("L37").Value = 15000
$ws.Range("M37").Value = -2618
$ws.Range("N37").Value = -15214

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2200.2632
$ws.Range("I55").Value = 1482.8889
$ws.Range("J55").Value = 2845.9
$ws.Range("K55").Value = 1482.8889
$ws.Range("L55").Value = 2845.9
$ws.Range("M55").Value = -1309.8889
$ws.Range("N55").Value = -3191.9

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2299.875
$ws.Range("I82").Value = 2520
$ws.Range("J82").Value = 1933
$ws.Range("K82").Value = 2520
$ws.Range("L82").Value = 1933
$ws.Range("M82").Value = -2159
$ws.Range("N82").Value = -2655

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2299.875
$ws.Range("I85").Value = 2520
$ws.Range("J85").Value = 1933
$ws.Range("K85").Value = 2520
$ws.Range("L85").Value = 1933
$ws.Range("M85").Value = -1272
$ws.Range("N85").Value = -4429

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3066.6667
$ws.Range("I93").Value = 3100
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 3100
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = -1852
$ws.Range("N93").Value = -5496

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6249.5
$ws.Range("I132").Value = 8000
$ws.Range("J132").Value = 4499
$ws.Range("K132").Value = 24000
$ws.Range("L132").Value = 13497
$ws.Range("M132").Value = -21470
$ws.Range("N132").Value = -18557

# WVR row 135
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").ClearContents()
$ws.Range("N135").ClearContents()

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5964.125
$ws.Range("I136").Value = 5062.5
$ws.Range("J136").Value = 6865.75
$ws.Range("K136").Value = 15187.5
$ws.Range("L136").Value = 20597.25
$ws.Range("M136").Value = -12637.5
$ws.Range("N136").Value = -25697.25
